$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, shifting everything down
$ws.Rows(1).Insert()

# Fill in the new header row (B1:I1) with the unit labels
$ws.Range("B1").Value = "U (PE-D)"
$ws.Range("C1").Value = "U (HG-S)"
$ws.Range("D1").Value = "U (HG-D)"
$ws.Range("E1").Value = "U (PC-S)"
$ws.Range("F1").Value = "U (PC-D)"
$ws.Range("G1").Value = "U (PVC-S)"
$ws.Range("H1").Value = "U (PVC-S)"
$ws.Range("I1").Value = "U (PVC-D)"
